# Legacy GSC export data refresh: the "Chart" sheet's oldest day
# (2025-10-06, row 2) drops off the rolling date window. Delete that
# row so every later day's data shifts up by one — matching a fresh
# export that now starts at 2025-10-07 and no longer includes the
# trailing 2025-12-23 row (the table shrinks from 80 to 79 rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")
$ws.Rows.Item(2).Delete()
